$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data values in row 2 (offset amplitude columns A,C,E,G,I,K,M) ---
$ws.Range("A2").Value = 4961272.0199999996
$ws.Range("C2").Value = 14866726.869999999
$ws.Range("E2").Value = 24774372.23
$ws.Range("G2").Value = 34680641.270000003
$ws.Range("I2").Value = 44587621.439999998
$ws.Range("K2").Value = 54494404.460000001
$ws.Range("M2").Value = 64401754.57

# --- Apply a medium box border + centered/wrapped alignment to those same cells ---
$a2 = $ws.Range("A2")
$a2.Borders.Weight = -4138
$a2.VerticalAlignment = -4108
$a2.WrapText = $true

$a2.Copy()
$cols = @("C","E","G","I","K","M")
foreach ($col in $cols) {
    $ws.Range($col + "2").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Header row (row 1) goes back to the default (unstyled) look ---
$ws.Range("A1:N1").Style = "Normal"

# --- Row heights grow slightly to fit the new bottom border look ---
$ws.Rows("1:2").RowHeight = 15.75

# --- Column E widens to fit the new large numbers ---
$ws.Columns("E").ColumnWidth = 10.736979166666666

# --- Update the active selection to M2 ---
$ws.Range("M2").Select()
